$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data was double-UTF8-encoded, turning the intended "±"
# (U+00B1) glyph into the two-character mojibake sequence "Â±"
# (U+00C2, U+00B1). Repair it back to the single correct character in
# every affected cell (columns B, C, D across rows 2-17).
$bad  = [string]([char]0x00C2) + [string]([char]0x00B1)
$good = [string]([char]0x00B1)

$cols = "B", "C", "D"

for ($row = 2; $row -le 17; $row++) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $cell = $ws.Range($addr)
        $val = $cell.Value2
        if ($val -ne $null -and $val.GetType().FullName -eq "System.String" -and $val.IndexOf($bad) -ge 0) {
            $cell.Value = $val.Replace($bad, $good)
        }
    }
}
